$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row 3: rename column header text ---
# H3 keeps the same text "RandomCardRate" (shared-string reindex only, no visible change)
$ws.Range("H3").Value = "RandomCardRate"
# I3 header text changes from "RandomCardRate2" to "RandomCardCatalog"
$ws.Range("I3").Value = "RandomCardCatalog"

# --- Rename table column (ListObject) "RandomCardRate2" -> "RandomCardCatalog" ---
$tbl = $ws.ListObjects.Item(1)
$col = $tbl.ListColumns.Item("RandomCardRate2")
$col.Name = "RandomCardCatalog"

# --- Rows 4-10: update H (rate) and I (catalog) values ---
$ws.Range("H4").Value = "0;840;150;10;10"
$ws.Range("I4").Value = "1;1;0"

$ws.Range("H5").Value = "0;840;150;10;10"
$ws.Range("I5").Value = "1;1;1"

$ws.Range("H6").Value = "0;840;150;10;10"
$ws.Range("I6").Value = "1;1;2"

$ws.Range("H7").Value = "0;840;150;10;10"
$ws.Range("I7").Value = "1;1;3"

$ws.Range("H8").Value = "0;840;150;10;10"
$ws.Range("I8").Value = "1;1;4"

$ws.Range("H9").Value = "0;840;150;10;10"
$ws.Range("I9").Value = "1;1;5"

$ws.Range("H10").Value = "0;840;150;10;10"
$ws.Range("I10").Value = "1;1;6"

# --- Rows 48-50: update H (rate) values and add new I (catalog) values ---
$ws.Range("H48").Value = "0;840;150;10;10"
$ws.Range("I48").Value = "5;0;0"

$ws.Range("H49").Value = "0;720;250;30;10"
$ws.Range("I49").Value = "5;0;0"

$ws.Range("H50").Value = "0;600;350;50;10"
$ws.Range("I50").Value = "5;0;0"

# --- sheetView: scroll position topLeftCell A37 -> A40, keep selection H50 ---
$ws.Activate()
$ws.Range("A40").Select()
$ws.Range("H50").Select()

$wb.Save()
